$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.660.57"
$ws.Range("E2").Value = "  -5.08%  "

$ws.Range("D3").Value = "3.386.35"
$ws.Range("E3").Value = "  -2.31%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -7.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "631.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.52%  "

$ws.Range("E7").Value = "  -6.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.392"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.86%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.949"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.30%  "

$ws.Range("D11").Value = "3.385.76"
$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("E12").Value = "  -7.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -12.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.41%  "

$ws.Range("D15").Value = "93.473.24"
$ws.Range("E15").Value = "  -5.12%  "

$ws.Range("D16").Value = "4.020.10"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000246"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -11.83%  "

$ws.Range("D19").Value = "3.389.16"
$ws.Range("E19").Value = "  -2.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.458"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -15.27%  "

$ws.Range("E24").Value = "  -8.53%  "

$ws.Range("E25").Value = "  -8.76%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "90.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.39%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("E31").Value = "  -14.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.133"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.987"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.174"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.535"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.64%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.34%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "529.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.90%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  -5.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.889"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.63%  "

$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.89%  "

$ws.Range("E47").Value = "  -4.06%  "

$ws.Range("E48").Value = "  -9.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.46%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.44%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.29%  "
